$p = $ppt.ActivePresentation
$p.Slides.Item(24).Delete()
